$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows above the old row 20 (the "view_data" / view_data_instructions
# block) to hold the new "data view modal" translation strings: dl_num_rows,
# frequency, daily, hourly, max.
$ws.Rows("20:24").Insert()

$ws.Range("A20").Value = "dl_num_rows"
$ws.Range("B20").Value = "advice about number of rows selected (modal)"
$ws.Range("C20").Value = "Number of rows that will be returned:"
$ws.Range("D20").Value = "Nombre de rangées à télécharger:"

$ws.Range("A21").Value = "frequency"
$ws.Range("B21").Value = "modal selection title"
$ws.Range("C21").Value = "Frequency:"
$ws.Range("D21").Value = "Fréquence:"

$ws.Range("A22").Value = "daily"
$ws.Range("B22").Value = "modal selection"
$ws.Range("C22").Value = "Daily"
$ws.Range("D22").Value = "Journalière"

$ws.Range("A23").Value = "hourly"
$ws.Range("B23").Value = "modal selection"
$ws.Range("C23").Value = "Hourly"
$ws.Range("D23").Value = "Horaire"

$ws.Range("A24").Value = "max"
$ws.Range("B24").Value = "modal selection"
$ws.Range("C24").Value = "Max"
$ws.Range("D24").Value = "Maximum"

# Update the view selection/scroll position to match the edited workbook.
$excel.Goto($ws.Range("A9"), $true)
$ws.Range("B24").Select()
